$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.585.94'
$ws.Range("E2").Value = '''  +0.35%  '
$ws.Range("D3").Value = '''1.923.76'
$ws.Range("E3").Value = '''  -0.69%  '
$ws.Range("E4").Value = '''  +0.07%  '
$ws.Range("D5").Value = '''247.36'
$ws.Range("E5").Value = '''  +2.83%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '''  +0.10%  '
$ws.Range("D7").Value = '''0.4730'
$ws.Range("E7").Value = '''  -0.33%  '
$ws.Range("D8").Value = '''0.2912'
$ws.Range("E8").Value = '''  +1.12%  '
$ws.Range("D9").Value = '''0.06834'
$ws.Range("E9").Value = '''  +2.85%  '
$ws.Range("D10").Value = '''105.87'
$ws.Range("E10").Value = '''  -0.90%  '
$ws.Range("D11").Value = '''18.40'
$ws.Range("D12").Value = '''1.931.63'
$ws.Range("E12").Value = '''  -0.16%  '
$ws.Range("D13").Value = '''0.07723'
$ws.Range("E13").Value = '''  +1.34%  '
$ws.Range("D14").Value = '''5.325'
$ws.Range("E14").Value = '''  +3.12%  '
$ws.Range("D15").Value = '''0.6714'
$ws.Range("E15").Value = '''  +1.02%  '
$ws.Range("D16").Value = '''293.05'
$ws.Range("E16").Value = '''  -5.17%  '
$ws.Range("D17").Value = '''30.622.11'
$ws.Range("E17").Value = '''  +0.42%  '
$ws.Range("D18").Value = '''0.000007643'
$ws.Range("E18").Value = '''  +0.79%  '
$ws.Range("D19").Value = '''5.577'
$ws.Range("E19").Value = '''  +5.14%  '
$ws.Range("D20").Value = '''0.9993'
$ws.Range("E20").Value = '''  -0.05%  '
$ws.Range("E21").Value = '''  -0.56%  '
$ws.Range("D22").Value = '''2.173.27'
$ws.Range("E22").Value = '''  -0.24%  '
$ws.Range("E23").Value = '''  +0.00%  '
$ws.Range("D24").Value = '''6.489'
$ws.Range("E24").Value = '''  +2.72%  '
$ws.Range("D25").Value = '''9.530'
$ws.Range("E25").Value = '''  +2.14%  '
$ws.Range("E26").Value = '''  -0.23%  '
$ws.Range("D27").Value = '''20.91'
$ws.Range("E27").Value = '''  +2.99%  '
$ws.Range("D28").Value = '''2.130'
$ws.Range("E28").Value = '''  +3.42%  '
$ws.Range("D29").Value = '''0.1071'
$ws.Range("E29").Value = '''  -3.42%  '
$ws.Range("D30").Value = '''1.405'
$ws.Range("E30").Value = '''  +2.53%  '
$ws.Range("D31").Value = '''4.195'
$ws.Range("E31").Value = '''  +2.01%  '
$ws.Range("D32").Value = '''4.060'
$ws.Range("E32").Value = '''  +2.97%  '
$ws.Range("D33").Value = '''0.05034'
$ws.Range("E33").Value = '''  +0.01%  '
$ws.Range("D34").Value = '''0.7350'
$ws.Range("E34").Value = '''  -1.09%  '
$ws.Range("D35").Value = '''1.147'
$ws.Range("E35").Value = '''  -0.91%  '
$ws.Range("D36").Value = '''0.02064'
$ws.Range("E36").Value = '''  +4.88%  '
$ws.Range("E37").Value = '''  +0.06%  '
$ws.Range("D38").Value = '''2.730'
$ws.Range("E38").Value = '''  -1.17%  '
$ws.Range("D39").Value = '''2.684'
$ws.Range("E39").Value = '''  -0.36%  '
$ws.Range("B40").Value = 'Quant'
$ws.Range("C40").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D40").Value = '''111.65'
$ws.Range("E40").Value = '''  +3.58%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '''111.65'
$ws.Range("E41").Value = '''  -0.18%  '
$ws.Range("D42").Value = '''0.4453'
$ws.Range("E42").Value = '''  +6.26%  '
$ws.Range("D43").Value = '''0.8731'
$ws.Range("E43").Value = '''  -0.98%  '
$ws.Range("D44").Value = '''5.896'
$ws.Range("E44").Value = '''  +1.41%  '
$ws.Range("D45").Value = '''1.0000'
$ws.Range("E45").Value = '''  +0.06%  '
$ws.Range("E46").Value = '''  -3.64%  '
$ws.Range("D47").Value = '''7.304'
$ws.Range("E47").Value = '''  +0.08%  '
$ws.Range("D48").Value = '''9.394'
$ws.Range("E48").Value = '''  +1.61%  '
$ws.Range("D49").Value = '''0.1253'
$ws.Range("E49").Value = '''  +3.21%  '
$ws.Range("D50").Value = '''47.89'
$ws.Range("E50").Value = '''  +12.67%  '
$ws.Range("D51").Value = '''35.21'
$ws.Range("E51").Value = '''  +0.78%  '
